$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "Fn15 / Cheque Services" row (old row 16); rows below shift up by one.
$ws.Rows.Item(16).Delete()

# Normalize function-name labels in column B to use underscores, and fix two
# description typos in column C (post-deletion row numbers).
$ws.Range("B2").Value = 'Verify_PIN'
$ws.Range("B3").Value = 'View_Account '
$ws.Range("B4").Value = 'Transfer_Money'
$ws.Range("B5").Value = 'Withdraw_cash'
$ws.Range("B6").Value = 'Apply_loan'
$ws.Range("B7").Value = 'Restock_cash'
$ws.Range("B9").Value = 'Cheque_Services'
$ws.Range("C9").Value = 'Cheque service is the functionality by which the customer may enquiries cheque status, whether it is paid, unpaid, stopped or returned. It also allows customer to stop cheque payment and to order cheque book to be delivered at home .The customer must be logged into Banking System.'
$ws.Range("B10").Value = 'Review_transactions'
$ws.Range("C10").Value = 'If the customer wants to display his/her payment history, review old transactions after withdraw amount or deposit etc. just he/she has to click on Bill Payment History, the system will display the transaction he/she done.'
$ws.Range("B11").Value = 'Limit_Cash'
$ws.Range("B13").Value = 'View_Account '
$ws.Range("B14").Value = 'Transfer_Funds '
$ws.Range("B15").Value = 'Pay_Bills '
$ws.Range("B19").Value = 'Pay_Registered_Payment'
$ws.Range("B20").Value = 'Open_Payment'
$ws.Range("B21").Value = 'Pay_Registration_Bill'
$ws.Range("B22").Value = 'Delete_registration_Bill '
$ws.Range("B24").Value = 'Transaction_Details'
$ws.Range("B25").Value = 'Deposit_Cheque'

# Restore the active selection to B11, matching the saved view state.
$ws.Range("B11").Select()
